$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.763.31'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '2.564.52'
$ws.Range('E3').Value = '  -3.15%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'515.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').Value = "'141.96"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.83%  '
$ws.Range('D8').Value = "'0.565"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').Value = '2.579.51'
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('D10').Value = "'6.59"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('D11').Value = "'0.100"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  -4.78%  '
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').Value = '3.018.43'
$ws.Range('E14').Value = '  -3.16%  '
$ws.Range('D15').Value = '57.801.78'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').Value = "'20.21"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.58%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = "'0.0000133"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.50%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.562.14'
$ws.Range('E18').Value = '  -3.49%  '
$ws.Range('D19').Value = "'337.77"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D21').Value = "'10.17"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('D22').Value = "'6.30"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = "'0.998"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = "'65.29"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('D25').Value = "'0.165"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('E27').Value = '  -6.00%  '
$ws.Range('D28').Value = '2.685.14'
$ws.Range('E28').Value = '  -2.99%  '
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('D30').Value = '0.0₃0740'
$ws.Range('E30').Value = '  -7.32%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').Value = "'6.25"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.53%  '
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = "'149.82"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = "'18.62"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('E36').Value = '  -4.10%  '
$ws.Range('E37').Value = '  -3.89%  '
$ws.Range('D38').Value = "'0.867"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.65%  '
$ws.Range('D39').Value = "'36.04"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = "'0.828"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'1.44"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('E42').Value = '  -3.26%  '
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = "'269.49"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').Value = "'10.66"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').Value = "'0.0949"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.32%  '
$ws.Range('E47').Value = '  -3.88%  '
$ws.Range('D48').Value = "'18.69"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.15%  '
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').Value = '1.976.90'
$ws.Range('E50').Value = '  -3.32%  '
$ws.Range('D51').Value = "'4.58"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.39%  '
